$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# There is no unemployment data for this series before 2003, so drop the
# 2000-2002 rows (rows 2:4) and let everything below shift up.
$ws.Rows("2:4").Delete()

# Restore the active selection to match what was left selected after the edit.
$ws.Range("A2:XFD4").Select()
